$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = -6.715099999999995
$ws.Range("A3").Value = -21.28310000000003
$ws.Range("D5").Value = -8.307499999999994
$ws.Range("A14").Value = -20.55009999999998
$ws.Range("A16").Value = -20.16469999999999
$ws.Range("D16").Value = -7.957700000000002
$ws.Range("A21").Value = -21.15670000000001
$ws.Range("A23").Value = -21.48290000000002
$ws.Range("A25").Value = -22.43340000000003
